# Diary update:
#  - appends a second sentence to the last existing entry
#    (06/11/2024) as its own run
#  - adds a new blank line, a new "08/11/2024" entry and a
#    trailing blank line
#
# We build the new content as literal OOXML paragraph fragments and
# splice them in with Range.InsertXML so that run boundaries (and the
# genuinely-empty <w:p/> separator paragraphs) come out exactly as
# authored, rather than being coalesced the way InsertAfter/TypeText
# would merge same-formatted runs together.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) Replace the last paragraph's content with the original sentence
#    plus the new trailing run, keeping the paragraph itself (and its
#    identity/attributes) intact.
$lastPara = $d.Paragraphs.Last
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$updatedParaXml = '<w:p ' + $wNs + '>' +
    '<w:r><w:t>Continue work on UML and code for bandit problem proof of concept</w:t></w:r>' +
    '<w:r><w:t>, did the explore and exploit only classes.</w:t></w:r>' +
    '</w:p>'
$lastRange.InsertXML($updatedParaXml)

# 2) Append the new blank line + new diary entry + trailing blank
#    line right after it.
$afterPara = $d.Paragraphs.Last
$insertPoint = $d.Range($afterPara.Range.End, $afterPara.Range.End)

$newBlocksXml = '<w:p ' + $wNs + '/>' +
    '<w:p ' + $wNs + '><w:r><w:t>08/11/2024</w:t></w:r></w:p>' +
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t xml:space="preserve">Continue work on UML and code for bandit problem proof of concept, did the </w:t></w:r>' +
        '<w:r><w:t>E-Greedy class</w:t></w:r>' +
        '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wNs + '/>'
$insertPoint.InsertXML($newBlocksXml)
